$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 58 (2025-06) with new retained_customers / retention_rate values
$ws.Range("B58").Value = 161
$ws.Range("D58").Value = 65.98360655737704

# Add new row 59 for 2025-07
$ws.Range("A59").Value = "2025-07"
$ws.Range("B59").Value = 1
$ws.Range("C59").Value = 227
$ws.Range("D59").Value = 0.4405286343612335
